# Update Ccl12-Ccr1 LR-pair sheet with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Ccr1", "ECs", "MuSCs", 1, 0.3333333333333333148296163, 0.2401406666666666966936106, 0.7204220000000000068141048, 0.0021927293339748929121291, 0.0022038442803109848690257, 3, 1, 30.346065666666660831651825, 91.0381969999999967058101902, 0.7437342022026641741305752, 0.7443504428118430382710358, 7.2873244399037773177951749, 65.5859199591340029655839317, 0.0016308078018501960104042, 0.0016404324659378289458234),
    @("FAPs", "Ccr1", "ECs", "Resolving-Mac", 1, 0.3333333333333333148296163, 0.2401406666666666966936106, 0.7204220000000000068141048, 0.0021927293339748929121291, 0.0022038442803109848690257, 1, 0.5, 0.1013394999999999990247801, 0.2026789999999999980495602, 0.0024836712940651791206581, 0.0016571528036595619196025, 0.0243357350896666706219218, 0.1460144105379999890370613, 0.000005446018902448099268, 0.0000036521067279464391354),
    @("FAPs", "Ccr1", "ECs", "Ccl12", 1, 0.3333333333333333148296163, 0.2401406666666666966936106, 0.7204220000000000068141048, 0.0021927293339748929121291, 0.0022038442803109848690257, 3, 1, 10.3548943333333305361065868, 31.064682999999998713747118, 0.2537821265032704776132277, 0.2539924043844973833294887, 2.4866312284695561984904089, 22.3796810562260013455215812, 0.0005564755132222484051985, 0.0005597597076452094167279),
    @("Inflammatory-Mac", "Ccr1", "ECs", "MuSCs", 1, 0.3333333333333333148296163, 0.643692000000000041914916, 1.9310760000000000147224455, 0.005877564804149374745601, 0.0059073581837392761043071, 3, 1, 30.346065666666660831651825, 91.0381969999999967058101902, 0.7437342022026641741305752, 0.7443504428118430382710358, 19.5335197011079984008574684, 175.8016773099720069239992881, 0.0043713459705084925857599, 0.0043971446799144948969218),
    @("Inflammatory-Mac", "Ccr1", "ECs", "Resolving-Mac", 1, 0.3333333333333333148296163, 0.643692000000000041914916, 1.9310760000000000147224455, 0.005877564804149374745601, 0.0059073581837392761043071, 1, 0.5, 0.1013394999999999990247801, 0.2026789999999999980495602, 0.0024836712940651791206581, 0.0016571528036595619196025, 0.0652314254339999982645892, 0.391388552604000017343111, 0.0000145979389830736301236, 0.0000097893951764048005412),
    @("Inflammatory-Mac", "Ccr1", "ECs", "Ccl12", 1, 0.3333333333333333148296163, 0.643692000000000041914916, 1.9310760000000000147224455, 0.005877564804149374745601, 0.0059073581837392761043071, 3, 1, 10.3548943333333305361065868, 31.064682999999998713747118, 0.2537821265032704776132277, 0.2539924043844973833294887, 6.6653626432120010036896929, 59.9882637889080072568503965, 0.0014916208946578069644007, 0.0015004241086483760883596),
    @("MuSCs", "Ccr1", "ECs", "MuSCs", 3, 1, 60.6610003333333267505622644, 181.9830010000000015679688659, 0.5538968334913180191492188, 0.5567045368793267501317246, 3, 1, 30.346065666666660831651825, 91.0381969999999967058101902, 0.7437342022026641741305752, 0.7443504428118430382710358, 1840.8226995210220593435224146, 16567.4042956892008078284561634, 0.4119520195592473177370607, 0.4143832685414888805652822),
    @("MuSCs", "Ccr1", "ECs", "Resolving-Mac", 3, 1, 60.6610003333333267505622644, 181.9830010000000015679688659, 0.5538968334913180191492188, 0.5567045368793267501317246, 1, 0.5, 0.1013394999999999990247801, 0.2026789999999999980495602, 0.0024836712940651791206581, 0.0016571528036595619196025, 6.1473554432798334090648495, 36.8841326596789969016754185, 0.0013756976652159870564557, 0.0009225444840995744143902),
    @("MuSCs", "Ccr1", "ECs", "Ccl12", 3, 1, 60.6610003333333267505622644, 181.9830010000000015679688659, 0.5538968334913180191492188, 0.5567045368793267501317246, 3, 1, 10.3548943333333305361065868, 31.064682999999998713747118, 0.2537821265032704776132277, 0.2539924043844973833294887, 628.1382486059648044829373248, 5653.2442374536831266595982015, 0.1405691162668545879377291, 0.1413987238537383062109143),
    @("Resolving-Mac", "Ccr1", "ECs", "MuSCs", 1, 0.5, 1.6570225000000000648014975, 3.314045000000000129602995, 0.015130306304387200824757, 0.0101380012242036198627515, 3, 1, 30.346065666666660831651825, 91.0381969999999967058101902, 0.7437342022026641741305752, 0.7443504428118430382710358, 50.2841135961441665358506725, 301.7046815768650276368134655, 0.0112529262883753599505665, 0.0075462257004629716347632),
    @("Resolving-Mac", "Ccr1", "ECs", "Resolving-Mac", 1, 0.5, 1.6570225000000000648014975, 3.314045000000000129602995, 0.015130306304387200824757, 0.0101380012242036198627515, 1, 0.5, 0.1013394999999999990247801, 0.2026789999999999980495602, 0.0024836712940651791206581, 0.0016571528036595619196025, 0.1679218316387499909669856, 0.6716873265550000748902448, 0.0000375787074386198967899, 0.0000168002171521930998587),
    @("Resolving-Mac", "Ccr1", "ECs", "Ccl12", 1, 0.5, 1.6570225000000000648014975, 3.314045000000000129602995, 0.015130306304387200824757, 0.0101380012242036198627515, 3, 1, 10.3548943333333305361065868, 31.064682999999998713747118, 0.2537821265032704776132277, 0.2539924043844973833294887, 17.1582928954558404655017512, 102.9497573727349930550190038, 0.0038398013085732250228299, 0.0025749753065884548909603),
    @("Ccl12", "Ccr1", "ECs", "MuSCs", 3, 1, 46.3149293333333389455219731, 138.9447879999999884148564888, 0.4229025660661706043086383, 0.4250462594324195264583466, 3, 1, 30.346065666666660831651825, 91.0381969999999967058101902, 0.7437342022026641741305752, 0.7443504428118430382710358, 1405.4758868963599525159224868, 12649.282982067239572643302381, 0.3145271025826827937521557, 0.3163833714240388839122886),
    @("Ccl12", "Ccr1", "ECs", "Resolving-Mac", 3, 1, 46.3149293333333389455219731, 138.9447879999999884148564888, 0.4229025660661706043086383, 0.4250462594324195264583466, 1, 0.5, 0.1013394999999999990247801, 0.2026789999999999980495602, 0.0024836712940651791206581, 0.0016571528036595619196025, 4.6935317811753334282798278, 28.1611906870519987933221273, 0.0010503509635250510029919, 0.0007043666005034435184662),
    @("Ccl12", "Ccr1", "ECs", "Ccl12", 3, 1, 46.3149293333333389455219731, 138.9447879999999884148564888, 0.4229025660661706043086383, 0.4250462594324195264583466, 3, 1, 10.3548943333333305361065868, 31.064682999999998713747118, 0.2537821265032704776132277, 0.2539924043844973833294887, 479.5861993024672074170666747, 4316.2757937222049804404377937, 0.1073251125199626010431331, 0.1079585214078770949441832)
)

$startRow = 2
$r = $startRow
foreach ($row in $data) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $col = $i + 1
        $val = $row[$i]
        $ws.Cells.Item($r, $col).Value2 = $val
    }
    $r = $r + 1
}

